# Scheduled runner update: refresh cached Universalis market-price snapshots
# (currentAveragePrice / LevePrice / LeveProfit columns, H:N) per-leve across sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 15833.667
$ws.Range("J13").Value = 15001
$ws.Range("L13").Value = 15001
$ws.Range("N13").Value = -15339

$ws.Range("H54").Value = 35559.8
$ws.Range("I54").Value = 36449.75
$ws.Range("K54").Value = 36449.75
$ws.Range("M54").Value = -35963.75

$ws.Range("H129").Value = 976.24194
$ws.Range("J129").Value = 1066.9615
$ws.Range("L129").Value = 3200.8845
$ws.Range("N129").Value = -13200.8845

$ws.Range("H132").Value = 1563.5294
$ws.Range("I132").Value = 1319.1936
$ws.Range("K132").Value = 3957.5808
$ws.Range("M132").Value = -1427.5808

$ws.Range("H134").Value = 111652
$ws.Range("J134").Value = 111652
$ws.Range("L134").Value = 111652
$ws.Range("N134").Value = -121792

$ws.Range("H138").Value = 3034556.8
$ws.Range("I138").Value = 6062450.5
$ws.Range("J138").Value = 6662.879
$ws.Range("K138").Value = 18187351.5
$ws.Range("L138").Value = 19988.637
$ws.Range("M138").Value = -18182211.5
$ws.Range("N138").Value = -30268.637

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14582.83
$ws.Range("I32").Value = 15441.358
$ws.Range("J32").Value = 6999.1665
$ws.Range("K32").Value = 15441.358
$ws.Range("L32").Value = 6999.1665
$ws.Range("M32").Value = -15154.358
$ws.Range("N32").Value = -7573.1665

$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H63").Value = 7251
$ws.Range("I63").Value = 6334.6665
$ws.Range("K63").Value = 6334.6665
$ws.Range("M63").Value = -5648.6665

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H66").Value = 7251
$ws.Range("I66").Value = 6334.6665
$ws.Range("K66").Value = 31673.3325
$ws.Range("M66").Value = -28241.3325

$ws.Range("H110").Value = 1696.238
$ws.Range("I110").Value = 1594.3572
$ws.Range("J110").Value = 1900
$ws.Range("K110").Value = 1594.3572
$ws.Range("L110").Value = 1900
$ws.Range("M110").Value = 450.6428000000001
$ws.Range("N110").Value = -5990

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 32003.941
$ws.Range("I107").Value = 38390.5
$ws.Range("J107").Value = 2200
$ws.Range("K107").Value = 38390.5
$ws.Range("L107").Value = 2200
$ws.Range("M107").Value = -36470.5
$ws.Range("N107").Value = -6040

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1188.15
$ws.Range("I58").Value = 1125.75
$ws.Range("K58").Value = 1125.75
$ws.Range("M58").Value = -922.75

$ws.Range("H136").Value = 1188.15
$ws.Range("I136").Value = 1125.75
$ws.Range("K136").Value = 3377.25
$ws.Range("M136").Value = -827.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1914.84
$ws.Range("J5").Value = 1053.8462
$ws.Range("L5").Value = 3161.5386
$ws.Range("N5").Value = -3385.5386

$ws.Range("H10").Value = 1700
$ws.Range("I10").Value = 50
$ws.Range("J10").Value = 5000
$ws.Range("K10").Value = 150
$ws.Range("L10").Value = 15000
$ws.Range("M10").Value = -11
$ws.Range("N10").Value = -15278

$ws.Range("H26").Value = 191.46666
$ws.Range("I26").Value = 114.333336
$ws.Range("J26").Value = 500
$ws.Range("K26").Value = 343.000008
$ws.Range("L26").Value = 1500
$ws.Range("M26").Value = -55.00000799999998
$ws.Range("N26").Value = -2076

$ws.Range("H113").Value = 883.6667
$ws.Range("I113").Value = 550
$ws.Range("J113").Value = 950.4
$ws.Range("K113").Value = 1650
$ws.Range("L113").Value = 2851.2
$ws.Range("M113").Value = 520
$ws.Range("N113").Value = -7191.2

$ws.Range("H117").Value = 52040.2
$ws.Range("J117").Value = 52040.2
$ws.Range("L117").Value = 156120.6
$ws.Range("N117").Value = -163004.6

$ws.Range("H121").Value = 34301.8
$ws.Range("I121").Value = 222
$ws.Range("J121").Value = 39544.848
$ws.Range("K121").Value = 666
$ws.Range("L121").Value = 118634.544
$ws.Range("M121").Value = 644
$ws.Range("N121").Value = -121254.544

$ws.Range("H131").Value = 881.64
$ws.Range("J131").Value = 887.89795
$ws.Range("L131").Value = 2663.69385
$ws.Range("N131").Value = -12743.69385

$ws.Range("H132").Value = 1778.2122
$ws.Range("I132").Value = 1006.8182
$ws.Range("J132").Value = 2163.9092
$ws.Range("K132").Value = 9061.363800000001
$ws.Range("L132").Value = 19475.1828
$ws.Range("M132").Value = -6531.363800000001
$ws.Range("N132").Value = -24535.1828

$ws.Range("H135").Value = 1914.84
$ws.Range("J135").Value = 1053.8462
$ws.Range("L135").Value = 9484.6158
$ws.Range("N135").Value = -14554.6158

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 116.666664
$ws.Range("I22").Value = 116.666664
$ws.Range("K22").Value = 116.666664
$ws.Range("M22").Value = 412.333336

$ws.Range("H51").Value = 43333.11
$ws.Range("J51").Value = 43333.11
$ws.Range("L51").Value = 43333.11
$ws.Range("N51").Value = -44351.11

$ws.Range("H132").Value = 1998.0392
$ws.Range("I132").Value = 1763
$ws.Range("J132").Value = 2468.1177
$ws.Range("K132").Value = 5289
$ws.Range("L132").Value = 7404.353099999999
$ws.Range("M132").Value = -2759
$ws.Range("N132").Value = -12464.3531

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3363.1667
$ws.Range("I7").Value = 2711.5
$ws.Range("J7").Value = 4666.5
$ws.Range("K7").Value = 2711.5
$ws.Range("L7").Value = 4666.5
$ws.Range("M7").Value = -2599.5
$ws.Range("N7").Value = -4890.5

$ws.Range("H97").Value = 25133.9
$ws.Range("J97").Value = 25133.9
$ws.Range("L97").Value = 25133.9
$ws.Range("N97").Value = -27115.9

$ws.Range("H126").Value = 3363.1667
$ws.Range("I126").Value = 2711.5
$ws.Range("J126").Value = 4666.5
$ws.Range("K126").Value = 8134.5
$ws.Range("L126").Value = 13999.5
$ws.Range("M126").Value = -5664.5
$ws.Range("N126").Value = -18939.5

$ws.Range("H134").Value = 76494.14
$ws.Range("J134").Value = 76494.14
$ws.Range("L134").Value = 76494.14
$ws.Range("N134").Value = -86634.14

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 17000.75
$ws.Range("I17").Value = 17000.75
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 17000.75
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -16828.75
$ws.Range("N17").ClearContents()
